$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (currently "Brown" stats row labeled via A2) becomes "ranking" with new values
$ws.Range("A2").Value = "ranking"
$ws.Range("B2").Value = 1.06615
$ws.Range("C2").Value = 2.0279
$ws.Range("D2").Value = 2.90595

# Row 3 stays "mean" but values change slightly
$ws.Range("A3").Value = "mean"
$ws.Range("B3").Value = 25.586605
$ws.Range("C3").Value = 26.928288
$ws.Range("D3").Value = 28.163975

# Row 4 becomes "normalised_mean" with new values (previously held "median" row data)
$ws.Range("A4").Value = "normalised_mean"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.052773423435958
$ws.Range("D4").Value = 1.101065611904747

# Row 5 (new) becomes "median" with the values that used to be in row 4
$ws.Range("A5").Value = "median"
$ws.Range("B5").Value = 25.35
$ws.Range("C5").Value = 26.9
$ws.Range("D5").Value = 28.4

# Row 6 (new) becomes "SEM" with values close to the original row 2 values
$ws.Range("A6").Value = "SEM"
$ws.Range("B6").Value = 0.4391213094404399
$ws.Range("C6").Value = 0.7028509443633343
$ws.Range("D6").Value = 0.5469988507134235
